$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D7").Value = "2016-03-04 08:32:06"
$wsZhCn.Range("G7").Value = "2016-03-04 08:32:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D7").Value = "2016-03-04 08:32:17"
$wsDeDe.Range("G7").Value = "2016-03-04 08:33:07"
